$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value "looks like a number" to Excel's type inference
# (single decimal point, digits only) need to be forced back to text so
# they keep their exact printed form (trailing zeros, etc.) instead of
# being silently parsed into a floating point number. We do this by
# temporarily switching the cell to a text format, writing the value,
# then resetting the cell style back to Normal (removing the style index
# again so the cell XML stays styleless, matching the original).
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "26.155.77"
$ws.Range("E2").Value = "  +3.73%  "

Set-TextValue "D3" "1.603.05"
$ws.Range("E3").Value = "  +3.39%  "

Set-TextValue "D5" "212.70"
$ws.Range("E5").Value = "  +3.00%  "

$ws.Range("E6").Value = "  -0.24%  "

Set-TextValue "D7" "0.487"
$ws.Range("E7").Value = "  +2.33%  "

$ws.Range("E8").Value = "  +3.11%  "

$ws.Range("E9").Value = "  +1.65%  "

Set-TextValue "D10" "18.01"
$ws.Range("E10").Value = "  +1.39%  "

$ws.Range("E11").Value = "  +5.04%  "

Set-TextValue "D12" "1.826.00"
$ws.Range("E12").Value = "  +3.40%  "

Set-TextValue "D13" "1.598.16"
$ws.Range("E13").Value = "  +3.33%  "

$ws.Range("E14").Value = "  +0.86%  "

$ws.Range("E15").Value = "  +1.80%  "

Set-TextValue "D16" "26.134.98"
$ws.Range("E16").Value = "  +3.79%  "

Set-TextValue "D17" "60.45"
$ws.Range("E17").Value = "  +3.20%  "

$ws.Range("E18").Value = "  +2.11%  "

$ws.Range("E19").Value = "  -0.21%  "

Set-TextValue "D20" "204.44"
$ws.Range("E20").Value = "  +10.36%  "

$ws.Range("E21").Value = "  +3.49%  "

$ws.Range("E22").Value = "  +1.01%  "

Set-TextValue "D23" "5.99"
$ws.Range("E23").Value = "  +2.84%  "

$ws.Range("E24").Value = "  +11.12%  "

Set-TextValue "D25" "141.79"
$ws.Range("E25").Value = "  +1.76%  "

$ws.Range("E26").Value = "  -0.19%  "

$ws.Range("E27").Value = "  -2.83%  "

Set-TextValue "D28" "15.19"

$ws.Range("E29").Value = "  +0.92%  "

$ws.Range("E30").Value = "  +2.08%  "

$ws.Range("E31").Value = "  +2.35%  "

$ws.Range("E32").Value = "  +3.69%  "

$ws.Range("E33").Value = "  +0.77%  "

$ws.Range("E34").Value = "  +2.14%  "

$ws.Range("E35").Value = "  +1.97%  "

Set-TextValue "D36" "0.0164"
$ws.Range("E36").Value = "  +10.42%  "

Set-TextValue "D37" "1.116.68"
$ws.Range("E37").Value = "  +2.89%  "

$ws.Range("E38").Value = "  +0.16%  "

Set-TextValue "D39" "2.31"
$ws.Range("E39").Value = "  +2.63%  "

$ws.Range("E40").Value = "  +3.34%  "

Set-TextValue "D41" "0.493"
$ws.Range("E41").Value = "  +0.02%  "

Set-TextValue "D42" "0.781"
$ws.Range("E42").Value = "  -3.08%  "

Set-TextValue "D43" "1.738.68"
$ws.Range("E43").Value = "  +3.42%  "

Set-TextValue "D44" "5.12"
$ws.Range("E44").Value = "  +1.62%  "

Set-TextValue "D45" "92.97"
$ws.Range("E45").Value = "  +0.49%  "

Set-TextValue "D46" "1.51"

Set-TextValue "D47" "53.43"

$ws.Range("E49").Value = "  +1.07%  "

$ws.Range("E50").Value = "  -0.07%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D51" ("0.0{0}0924" -f [char]0x2087)
$ws.Range("E51").Value = "  -17.17%  "
